$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Kane Williamson ", "57", "39", "7", "0", "146.15"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Sunrisers Hyderabad", "Kolkata Knight Riders", "Kane Williamson ", "29", "19", "4", "1", "152.63"),
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Kane Williamson ", "11", "10", "1", "0", "110.00"),
    @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Kane Williamson ", "9", "13", "1", "0", "69.23"),
    @(" Abu Dhabi", " September 29 2020", "Sunrisers won by 15 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Kane Williamson ", "41", "26", "5", "0", "157.69"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Sunrisers Hyderabad", "Mumbai Indians", "Kane Williamson ", "3", "5", "0", "0", "60.00"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Kane Williamson ", "50", "44", "2", "2", "113.63"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Kane Williamson ", "8", "14", "0", "0", "57.14"),
    @(" Abu Dhabi", " November 08 2020", "Capitals won by 17 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Kane Williamson ", "67", "45", "5", "4", "148.88"),
    @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", "Kane Williamson ", "22", "12", "0", "2", "183.33"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Kane Williamson ", "20", "10", "1", "1", "200.00")
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
